# SMD_Inventory update:
#  - Bump quantities for three existing parts.
#  - Add a new inventory line for a Texas Instruments linear regulator
#    (TPS7A8500RGRT) in a 20VQFN package.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantity corrections on existing rows ---
$ws.Cells.Item(33, 6).Value = 28
$ws.Cells.Item(35, 6).Value = 3
$ws.Cells.Item(69, 6).Value = 6

# --- New row 102: IC / Linear Voltage Regulator / TPS7A8500RGRT ---
$newRow = 102

$cA = $ws.Cells.Item($newRow, 1)
$cA.Value = "IC"
$cA.NumberFormat = "@"

$cB = $ws.Cells.Item($newRow, 2)
$cB.Value = "Linear Voltage Regulator"
$cB.NumberFormat = "@"

$cC = $ws.Cells.Item($newRow, 3)
$cC.Value = "TPS7A8500RGRT"
$cC.NumberFormat = "@"

$cD = $ws.Cells.Item($newRow, 4)
$cD.Value = "Texas Instruments"
$cD.NumberFormat = "@"

$cE = $ws.Cells.Item($newRow, 5)
$cE.Value = "20VQFN"
$cE.NumberFormat = "@"

$cF = $ws.Cells.Item($newRow, 6)
$cF.Value = 6

$cG = $ws.Cells.Item($newRow, 7)
$cG.Value = "IC REG LINEAR POS ADJ 4A 20VQFN"
$cG.NumberFormat = "@"
